# Commit message: Replaced all instances of the term "squad" with "unit",
# and "unit" with "model" - reflecting the official 40k ruleset language.
# Concretely (per the canonical diff) this renames the "Templar Units"
# worksheet to "Templar Models", and the active/selected tab moves from
# "Templar Melee Weapons" to the renamed sheet, with its selection moved
# to G27.

$wb = $excel.ActiveWorkbook

# Rename "Templar Units" -> "Templar Models"
$wsModels = $wb.Worksheets.Item("Templar Units")
$wsModels.Name = "Templar Models"

# Move the active selection/tab onto the renamed sheet, selecting G27
# (this also clears tabSelected from whatever sheet was previously active).
$wsModels.Activate()
$null = $wsModels.Range("G27").Select()
